$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 410
$ws.Range("I12").Value = 410
$ws.Range("K12").Value = 410
$ws.Range("M12").Value = -240
$ws.Range("H62").Value = 5333.3335
$ws.Range("J62").Value = 5000
$ws.Range("L62").Value = 5000
$ws.Range("N62").Value = -6248
$ws.Range("H65").Value = 5333.3335
$ws.Range("J65").Value = 5000
$ws.Range("L65").Value = 25000
$ws.Range("N65").Value = -31240
$ws.Range("H86").Value = 43484100
$ws.Range("I86").Value = 62502148
$ws.Range("K86").Value = 62502148
$ws.Range("M86").Value = -62501025
$ws.Range("H89").Value = 43484100
$ws.Range("I89").Value = 62502148
$ws.Range("K89").Value = 312510740
$ws.Range("M89").Value = -312505124
$ws.Range("H138").Value = 3778.6667
$ws.Range("I138").Value = 3375.12
$ws.Range("J138").Value = 4166.6924
$ws.Range("K138").Value = 10125.36
$ws.Range("L138").Value = 12500.0772
$ws.Range("M138").Value = -4985.360000000001
$ws.Range("N138").Value = -22780.0772

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 50000
$ws.Range("I19").Value = 50000
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 50000
$ws.Range("L19").Value = 0
$ws.Range("N19").Value = 0
$ws.Range("M19").Value = -49771
$ws.Range("H32").Value = 2628.1167
$ws.Range("I32").Value = 2528.8447
$ws.Range("J32").Value = 5507
$ws.Range("K32").Value = 2528.8447
$ws.Range("L32").Value = 5507
$ws.Range("M32").Value = -2241.8447
$ws.Range("N32").Value = -6081
$ws.Range("H61").Value = 1832.3925
$ws.Range("I61").Value = 1134.791
$ws.Range("K61").Value = 1134.791
$ws.Range("M61").Value = -922.7909999999999
$ws.Range("H74").Value = 104971.39
$ws.Range("I74").Value = 130936.35
$ws.Range("K74").Value = 130936.35
$ws.Range("M74").Value = -130062.35
$ws.Range("H77").Value = 104971.39
$ws.Range("I77").Value = 130936.35
$ws.Range("K77").Value = 654681.75
$ws.Range("M77").Value = -650313.75
$ws.Range("H110").Value = 3214.6667
$ws.Range("J110").Value = 8993
$ws.Range("L110").Value = 8993
$ws.Range("N110").Value = -13083
$ws.Range("H122").Value = 4426.8237
$ws.Range("I122").Value = 4585.9556
$ws.Range("K122").Value = 13757.8668
$ws.Range("M122").Value = -11307.8668
$ws.Range("H132").Value = 5546234
$ws.Range("I132").Value = 1283904.4
$ws.Range("J132").Value = 23813360
$ws.Range("K132").Value = 3851713.2
$ws.Range("L132").Value = 71440080
$ws.Range("M132").Value = -3849183.2
$ws.Range("N132").Value = -71445140
$ws.Range("H136").Value = 1832.3925
$ws.Range("I136").Value = 1134.791
$ws.Range("K136").Value = 3404.373
$ws.Range("M136").Value = -854.3729999999996

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 9
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()
$ws.Range("H20").Value = 21935246
$ws.Range("I20").Value = 28741322
$ws.Range("J20").Value = 4556.3335
$ws.Range("K20").Value = 28741322
$ws.Range("L20").Value = 4556.3335
$ws.Range("M20").Value = -28741075
$ws.Range("N20").Value = -5050.3335
$ws.Range("H99").Value = 81254.69500000001
$ws.Range("I99").Value = 169337.5
$ws.Range("J99").Value = 5755.143
$ws.Range("K99").Value = 169337.5
$ws.Range("L99").Value = 5755.143
$ws.Range("M99").Value = -167839.5
$ws.Range("N99").Value = -8751.143
$ws.Range("H107").Value = 2080179.1
$ws.Range("I107").Value = 2849988.2
$ws.Range("J107").Value = 1694.4
$ws.Range("K107").Value = 2849988.2
$ws.Range("L107").Value = 1694.4
$ws.Range("M107").Value = -2848068.2
$ws.Range("N107").Value = -5534.4
$ws.Range("H134").Value = 2472.3674
$ws.Range("I134").Value = 2103.7026
$ws.Range("K134").Value = 6311.1078
$ws.Range("M134").Value = -3776.1078
$ws.Range("H140").Value = 58790.832
$ws.Range("J140").Value = 58790.832
$ws.Range("L140").Value = 58790.832
$ws.Range("N140").Value = -69150.83199999999

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1572.5834
$ws.Range("I22").Value = 1513.4445
$ws.Range("J22").Value = 1750
$ws.Range("K22").Value = 1513.4445
$ws.Range("L22").Value = 1750
$ws.Range("M22").Value = -1163.4445
$ws.Range("N22").Value = -2450
$ws.Range("H29").Value = 20000
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 20000
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 20000
$ws.Range("M29").Value = 20000
$ws.Range("N29").Value = -20586
$ws.Range("H31").Value = 2504866.8
$ws.Range("I31").Value = 2800.7666
$ws.Range("K31").Value = 2800.7666
$ws.Range("M31").Value = -2505.7666
$ws.Range("H34").Value = 2504866.8
$ws.Range("I34").Value = 2800.7666
$ws.Range("K34").Value = 2800.7666
$ws.Range("M34").Value = -2598.7666
$ws.Range("H58").Value = 1658.2041
$ws.Range("I58").Value = 1093.5526
$ws.Range("K58").Value = 1093.5526
$ws.Range("M58").Value = -890.5526
$ws.Range("H134").Value = 1825.5079
$ws.Range("I134").Value = 1785.8948
$ws.Range("J134").Value = 2201.8333
$ws.Range("K134").Value = 5357.6844
$ws.Range("L134").Value = 6605.499899999999
$ws.Range("M134").Value = -2822.6844
$ws.Range("N134").Value = -11675.4999
$ws.Range("H136").Value = 1658.2041
$ws.Range("I136").Value = 1093.5526
$ws.Range("K136").Value = 3280.6578
$ws.Range("M136").Value = -730.6578

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 893.2
$ws.Range("J5").Value = 948.0833
$ws.Range("L5").Value = 2844.2499
$ws.Range("N5").Value = -3068.2499
$ws.Range("H133").Value = 4605.636
$ws.Range("I133").Value = 4310
$ws.Range("J133").Value = 4716.5
$ws.Range("K133").Value = 12930
$ws.Range("L133").Value = 14149.5
$ws.Range("M133").Value = -7870
$ws.Range("N133").Value = -24269.5
$ws.Range("H135").Value = 893.2
$ws.Range("J135").Value = 948.0833
$ws.Range("L135").Value = 8532.7497
$ws.Range("N135").Value = -13602.7497

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 39959.5
$ws.Range("J43").Value = 39959.5
$ws.Range("L43").Value = 39959.5
$ws.Range("N43").Value = -40261.5
$ws.Range("H70").Value = 29539942
$ws.Range("I70").Value = 50203350
$ws.Range("K70").Value = 50203350
$ws.Range("M70").Value = -50203080
$ws.Range("H73").Value = 29539942
$ws.Range("I73").Value = 50203350
$ws.Range("K73").Value = 50203350
$ws.Range("M73").Value = -50202414
$ws.Range("H113").Value = 4372.5
$ws.Range("J113").Value = 4360
$ws.Range("L113").Value = 4360
$ws.Range("N113").Value = -8700
$ws.Range("H122").Value = 4087.3928
$ws.Range("I122").Value = 2873.6365
$ws.Range("J122").Value = 8537.833000000001
$ws.Range("K122").Value = 8620.9095
$ws.Range("L122").Value = 25613.499
$ws.Range("M122").Value = -6170.9095
$ws.Range("N122").Value = -30513.499
$ws.Range("H126").Value = 49384980
$ws.Range("I126").Value = 2548
$ws.Range("K126").Value = 7644
$ws.Range("M126").Value = -5174
$ws.Range("H132").Value = 1774.641
$ws.Range("I132").Value = 1424.5
$ws.Range("K132").Value = 4273.5
$ws.Range("M132").Value = -1743.5
$ws.Range("H134").Value = 54113.145
$ws.Range("J134").Value = 54113.145
$ws.Range("L134").Value = 162339.435
$ws.Range("N134").Value = -167409.435

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("H43").Value = 20035.715
$ws.Range("I43").Value = 19208.334
$ws.Range("J43").Value = 25000
$ws.Range("K43").Value = 19208.334
$ws.Range("L43").Value = 25000
$ws.Range("M43").Value = -19015.334
$ws.Range("N43").Value = -25386
$ws.Range("H61").Value = 6844.4546
$ws.Range("I61").Value = 1554.7142
$ws.Range("K61").Value = 1554.7142
$ws.Range("M61").Value = -1352.7142
$ws.Range("H113").Value = 6844.4546
$ws.Range("I113").Value = 1554.7142
$ws.Range("K113").Value = 1554.7142
$ws.Range("M113").Value = 615.2858000000001
$ws.Range("H122").Value = 998.3333
$ws.Range("I122").Value = 333.33334
$ws.Range("J122").Value = 1663.3334
$ws.Range("K122").Value = 1000.00002
$ws.Range("L122").Value = 4990.0002
$ws.Range("M122").Value = 1449.99998
$ws.Range("N122").Value = -9890.0002
$ws.Range("H136").Value = 3299.2354
$ws.Range("I136").Value = 2829.4348
$ws.Range("K136").Value = 8488.304400000001
$ws.Range("M136").Value = -5938.304400000001

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2179.742
$ws.Range("I132").Value = 1557.0577
$ws.Range("J132").Value = 5417.7
$ws.Range("K132").Value = 4671.1731
$ws.Range("L132").Value = 16253.1
$ws.Range("M132").Value = -2141.1731
$ws.Range("N132").Value = -21313.1
